$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Enhancement")

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "make 6 keywords result instead of 5 for better display"

$ws.Range("A5").Select()
